# NBN FTTP tab working
# Adds "target" / "categories" / "subCategories" columns to the header row,
# adds two new data rows ("tabs", "btn-selection"), drops the leftover
# empty-but-styled G16/G17 cells, widens column P, and moves the
# selection/view to the newly-added area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New header cells on row 11 (O11:Q11) ---
$ws.Range("O11").Value = "target"
$ws.Range("P11").Value = "categories"
$ws.Range("Q11").Value = "subCategories"

# --- Drop the stray empty styled cells at G16 / G17 ---
# Clear() removes both content AND formatting so the <c> element itself
# disappears from the saved XML (ClearContents alone would keep s="1").
$ws.Range("G16").Clear()
$ws.Range("G17").Clear()

# --- New row 20: "tabs" ---
$ws.Range("E20").Value = "tabs"
$ws.Range("O20").Value = "y"
$ws.Range("P20").Value = "y"

# --- New row 21: "btn-selection" ---
$ws.Range("E21").Value = "btn-selection"
$ws.Range("O21").Value = "y"
$ws.Range("Q21").Value = "y"

# --- Column P (16) gets a custom width ---
$ws.Columns.Item(16).ColumnWidth = 13.5

# --- View / selection moves onto the new data ---
$ws.Range("O21").Select()
